# [Kadastro App] Yeni kayit eklendi: 1
# Adds the new "Kayitlar" (records) row to both the master "Kayitlar" sheet
# and the per-district "Anamur" sheet, matching the row that was recorded
# there (Birim = Anamur).

$wb = $excel.ActiveWorkbook

$values = @("1", "2025-08-18", "Anamur", "2", "2", "AİLE KONUTU", "EMİNE ALANLI KIRCILI (K.Mühendisi), KAYHAN KARTPAK (K.Teknisyeni)")

foreach ($sheetName in @("Kayitlar", "Anamur")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Stage the new row far below the real data as text (so every value -
    # including the numeric-looking ones like "1" and "2" - lands as a text
    # cell, matching the existing numberStoredAsText column convention)
    # without leaving a lingering quote-prefix/text style on the final cells.
    $scratch = $ws.Range("A1000:G1000")
    $scratch.NumberFormat = "@"
    for ($i = 0; $i -lt $values.Length; $i++) {
        $scratch.Cells.Item(1, $i + 1).Value = $values[$i]
    }

    # Copy as values only into row 2 so the destination cells pick up the
    # plain "text" cell type without inheriting the scratch area's style.
    $scratch.Copy()
    $ws.Range("A2").PasteSpecial(-4163)

    # Remove the scratch row entirely so it leaves no trace.
    $ws.Rows.Item(1000).Delete()
}
